$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.503.52"
$ws.Range("E2").Value = "  +2.92%  "

# Row 3
$ws.Range("D3").Value = "1.839.28"
$ws.Range("E3").Value = "  +1.77%  "

# Row 4
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").Value = "231.64"
$ws.Range("E5").Value = "  +2.98%  "

# Row 6
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  +1.63%  "

# Row 7
$ws.Range("E7").Value = "  +0.31%  "

# Row 8
$ws.Range("D8").Value = "43.67"
$ws.Range("E8").Value = "  +12.23%  "

# Row 9
$ws.Range("E9").Value = "  +7.60%  "

# Row 10
$ws.Range("E10").Value = "  +5.04%  "

# Row 11
$ws.Range("E11").Value = "  +2.39%  "

# Row 12
$ws.Range("D12").Value = "2.104.62"
$ws.Range("E12").Value = "  +1.75%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.675"
$ws.Range("E13").Value = "  +6.80%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.835.39"
$ws.Range("E14").Value = "  +1.55%  "

# Row 15
$ws.Range("E15").Value = "  +1.38%  "

# Row 16
$ws.Range("E16").Value = "  +7.37%  "

# Row 17
$ws.Range("D17").Value = "35.457.98"
$ws.Range("E17").Value = "  +2.84%  "

# Row 18
$ws.Range("D18").Value = "70.14"
$ws.Range("E18").Value = "  +2.80%  "

# Row 19
$ws.Range("E19").Value = "  +4.26%  "

# Row 20
$ws.Range("D20").Value = "244.17"
$ws.Range("E20").Value = "  +1.45%  "

# Row 21
$ws.Range("D21").Value = "12.08"
$ws.Range("E21").Value = "  +8.16%  "

# Row 22
$ws.Range("D22").Value = "4.71"
$ws.Range("E22").Value = "  +14.61%  "

# Row 23
$ws.Range("E23").Value = "  +0.34%  "

# Row 24
$ws.Range("E24").Value = "  +0.84%  "

# Row 25
$ws.Range("D25").Value = "171.66"
$ws.Range("E25").Value = "  +0.35%  "

# Row 26
$ws.Range("D26").Value = "7.93"
$ws.Range("E26").Value = "  +3.09%  "

# Row 27
$ws.Range("E27").Value = "  +0.48%  "

# Row 28
$ws.Range("E28").Value = "  -0.98%  "

# Row 29
$ws.Range("D29").Value = "1.60"
$ws.Range("E29").Value = "  +30.78%  "

# Row 30
$ws.Range("E30").Value = "  +0.34%  "

# Row 31
$ws.Range("D31").Value = "3.313.10"
$ws.Range("E31").Value = "  +36.36%  "

# Row 32
$ws.Range("D32").Value = "0.0553"
$ws.Range("E32").Value = "  +7.35%  "

# Row 33
$ws.Range("E33").Value = "  +6.12%  "

# Row 34
$ws.Range("E34").Value = "  +4.69%  "

# Row 35
$ws.Range("E35").Value = "  +1.41%  "

# Row 36
$ws.Range("D36").Value = "95.40"
$ws.Range("E36").Value = "  +15.82%  "

# Row 37
$ws.Range("D37").Value = "0.690"
$ws.Range("E37").Value = "  +7.35%  "

# Row 38
$ws.Range("E38").Value = "  +5.93%  "

# Row 39
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.350.65"
$ws.Range("E39").Value = "  +3.24%  "

# Row 40
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "15.56"
$ws.Range("E40").Value = "  +11.31%  "

# Row 41
$ws.Range("E41").Value = "  +5.74%  "

# Row 42
$ws.Range("E42").Value = "  +4.77%  "

# Row 43
$ws.Range("E43").Value = "  +6.08%  "

# Row 44
$ws.Range("E44").Value = "  +4.48%  "

# Row 45
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  +1.07%  "

# Row 46
$ws.Range("E46").Value = "  +0.50%  "

# Row 47
$ws.Range("D47").Value = "6.26"
$ws.Range("E47").Value = "  +8.14%  "

# Row 48
$ws.Range("D48").Value = "0.0520"
$ws.Range("E48").Value = "  +1.11%  "

# Row 49
$ws.Range("D49").Value = "2.007.38"
$ws.Range("E49").Value = "  +1.99%  "

# Row 50
$ws.Range("E50").Value = "  +0.35%  "

# Row 51
$ws.Range("D51").Value = "103.34"
$ws.Range("E51").Value = "  +0.55%  "
